# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates the record_atd (column C) and average_simulation_TD-related (column D)
# values in the time-to-discovery simulation results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new value for columns C and D (both set to the same corrected value)
$updates = @{
    3  = 157
    5  = 54
    7  = 55
    9  = 185
    11 = 130
    13 = 129
    15 = 105
    17 = 63
    19 = 149
    21 = 150
    23 = 143
    25 = 85
    27 = 208
    29 = 183
    30 = 142
    32 = 299
    34 = 106
    36 = 175
    38 = 69
    40 = 147
    42 = 59
    44 = 112
    46 = 115
    48 = 193
    50 = 89
    52 = 84
    54 = 168
    56 = 691
    58 = 145
    60 = 88
    62 = 104
    64 = 87
    66 = 95
    68 = 156
    70 = 140
    72 = 134
    74 = 67
    76 = 99
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range("C$row").Value = $value
    $ws.Range("D$row").Value = $value
}

# Row 77 holds the average of column C only (no corresponding D value)
$ws.Range("C77").Value = 141.0526315789474
